$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 352. This shifts the existing rows 352-415
# down to 353-416 (carrying their values/formatting with them), exactly
# like Excel's native "Insert Row" behaviour.
$ws.Rows("352:352").Insert()

# Populate the newly inserted row 352 with the new weekly data point.
$ws.Range("A352").Value = 4
$ws.Range("B352").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C352").Value = "Los Lagos"
$ws.Range("D352").Value = 44964
$ws.Range("E352").Value = 10
$ws.Range("F352").Value = "Fruta"
$ws.Range("G352").Value = 100104
$ws.Range("H352").Value = "Frutos de pepita"
$ws.Range("I352").Value = 100104005
$ws.Range("J352").Value = "Pera"
$ws.Range("K352").Value = "Packham's Triumph"
$ws.Range("L352").Value = "Primera"
$ws.Range("M352").Value = 300
$ws.Range("N352").Value = 20000
$ws.Range("O352").Value = 21000
$ws.Range("P352").Value = 20500
$ws.Range("Q352").Value = "$/caja 15 kilos empedrada"
$ws.Range("R352").Value = "Región de O'Higgins"
$ws.Range("S352").Value = 1367
$ws.Range("T352").Value = 15
